$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.939.55"
$ws.Range("E2").Value = "  +2.81%  "

$ws.Range("D3").Value = "3.785.00"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.50%  "

$ws.Range("D7").Value = "3.784.29"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.88%  "

$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("E13").Value = "  +7.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.90%  "

$ws.Range("D15").Value = "4.421.57"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").Value = "3.784.30"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").Value = "70.924.29"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("E18").Value = "  +1.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +18.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "483.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.71%  "

$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.87%  "

$ws.Range("D29").Value = "3.935.06"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.92%  "

$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("E33").Value = "  +7.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.02%  "

$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.45%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "3.734.01"
$ws.Range("E38").Value = "  +0.43%  "

$ws.Range("E39").Value = "  +2.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000329"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +24.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.969"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "160.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "49.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("E51").Value = "  +1.89%  "
